# clean code and add comments
# Apply three small textual corrections to the results document:
#   1. Update the report date from Feb 24 to Feb 25, 2025
#   2. Rename "column specifications" to "variable specifications" in the
#      query info paragraph (matches the "Variable name" table header)
#   3. Update the processing time reported in the summary footer

$d = $word.ActiveDocument

# 1) Fix the report date shown under the title (Heading1 paragraph)
$d.Content.Find.Execute(
    "February 24, 2025", $true, $false, $false, $false, $false,
    $true, 1, $false, "February 25, 2025", 2
)

# 2) Correct terminology: "column specifications" -> "variable specifications"
$d.Content.Find.Execute(
    "The following query is run for each of the column specifications listed below:",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "The following query is run for each of the variable specifications listed below:",
    2
)

# 3) Update the elapsed processing time in the run summary heading
$d.Content.Find.Execute(
    "1 documents (44 total pages) processed in 4.85 seconds", $true, $false, $false, $false, $false,
    $true, 1, $false, "1 documents (44 total pages) processed in 6.48 seconds", 2
)
